$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.784.51'
$ws.Range('E2').Value = '  +3.11%  '
$ws.Range('D3').Value = '2.543.53'
$ws.Range('E3').Value = '  +5.67%  '
$ws.Range('E4').Value = '  -0.02%  '
$ws.Range('D5').Value = "'" + '573.60'
$ws.Range('E5').Value = '  +2.70%  '
$ws.Range('D6').Value = "'" + '148.42'
$ws.Range('E6').Value = '  +7.66%  '
$ws.Range('D8').Value = "'" + '0.590'
$ws.Range('E8').Value = '  +0.90%  '
$ws.Range('D9').Value = '2.543.23'
$ws.Range('E9').Value = '  +5.78%  '
$ws.Range('E10').Value = '  +2.69%  '
$ws.Range('D11').Value = "'" + '5.75'
$ws.Range('E11').Value = '  +0.24%  '
$ws.Range('E12').Value = '  +2.04%  '
$ws.Range('E13').Value = '  +3.05%  '
$ws.Range('E14').Value = '  +9.27%  '
$ws.Range('D15').Value = '2.999.55'
$ws.Range('E15').Value = '  +5.70%  '
$ws.Range('D16').Value = '63.583.42'
$ws.Range('E16').Value = '  +2.82%  '
$ws.Range('E17').Value = '  +3.67%  '
$ws.Range('D18').Value = '2.544.78'
$ws.Range('E18').Value = '  +5.41%  '
$ws.Range('D19').Value = "'" + '11.54'
$ws.Range('E19').Value = '  +4.27%  '
$ws.Range('D20').Value = "'" + '341.41'
$ws.Range('E20').Value = '  -0.34%  '
$ws.Range('E21').Value = '  +3.38%  '
$ws.Range('D22').Value = "'" + '6.87'
$ws.Range('E22').Value = '  +0.51%  '
$ws.Range('E23').Value = '  -0.07%  '
$ws.Range('D24').Value = "'" + '66.19'
$ws.Range('E24').Value = '  +1.98%  '
$ws.Range('D25').Value = "'" + '0.171'
$ws.Range('E25').Value = '  -0.80%  '
$ws.Range('E26').Value = '  +4.55%  '
$ws.Range('E27').Value = '  -0.03%  '
$ws.Range('D28').Value = "'" + '8.37'
$ws.Range('E28').Value = '  +0.92%  '
$ws.Range('D29').Value = "'" + '1.42'
$ws.Range('E29').Value = '  +3.16%  '
$ws.Range('D30').Value = '0.0₃0834'
$ws.Range('E30').Value = '  +7.85%  '
$ws.Range('D31').Value = "'" + '6.91'
$ws.Range('E31').Value = '  +9.15%  '
$ws.Range('E32').Value = '  +4.42%  '
$ws.Range('D33').Value = "'" + '177.77'
$ws.Range('E33').Value = '  +4.06%  '
$ws.Range('E34').Value = '  +14.33%  '
$ws.Range('D35').Value = "'" + '421.76'
$ws.Range('E35').Value = '  +12.70%  '
$ws.Range('D36').Value = "'" + '0.406'
$ws.Range('E36').Value = '  +2.52%  '
$ws.Range('D37').Value = "'" + '19.08'
$ws.Range('E37').Value = '  +3.07%  '
$ws.Range('E38').Value = '  -1.92%  '
$ws.Range('E40').Value = '  +6.39%  '
$ws.Range('D41').Value = "'" + '1.00'
$ws.Range('E41').Value = '  -0.01%  '
$ws.Range('D42').Value = "'" + '40.70'
$ws.Range('E42').Value = '  +4.23%  '
$ws.Range('D43').Value = "'" + '153.17'
$ws.Range('E43').Value = '  +6.23%  '
$ws.Range('E44').Value = '  +4.02%  '
$ws.Range('D45').Value = "'" + '20.88'
$ws.Range('E45').Value = '  +1.49%  '
$ws.Range('E46').Value = '  +4.47%  '
$ws.Range('D47').Value = "'" + '0.0533'
$ws.Range('E47').Value = '  +1.38%  '
$ws.Range('D48').Value = "'" + '0.0968'
$ws.Range('E48').Value = '  +1.25%  '
$ws.Range('B49').Value = 'VeChain'
$ws.Range('C49').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D49').Value = "'" + '0.0232'
$ws.Range('E49').Value = '  +5.54%  '
$ws.Range('B50').Value = 'BabyDogeCoin'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D50').Value = '0.0₆0239'
$ws.Range('E50').Value = '  +8.79%  '
$ws.Range('B51').Value = 'EnergySwap'
$ws.Range('C51').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D51').Value = "'" + '18.75'
$ws.Range('E51').Value = '  +5.04%  '
